# Updated cryptos list -- applies per-cell Price (D) and Volume(1h) (E) changes.
# NumberFormat "@" forces text interpretation so numeric-looking strings (e.g. "223.61")
# are not auto-converted to numbers by Excel's smart input; Style reset afterwards avoids
# leaving a stray style/number-format behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = '@'
$cell.Value = '34.536.95'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.38%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.786.41'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.92%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.12%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '223.61'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.42%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.559'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.29%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '32.95'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +8.25%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.281'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.48%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0678'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +3.06%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0937'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.51%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.043.86'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.03%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.07'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +11.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.792.79'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.30%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.21%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '34.527.43'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.54%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.83%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '68.49'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.37%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '253.43'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.08%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0774'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +5.35%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.01%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '10.41'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.90%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.23'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.41%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.14'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.40%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.13%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.53%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.73%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.14%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.10%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.75%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.98%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.38%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.58'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.50%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.85'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +3.59%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.443.86'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.49%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.42%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0189'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.51%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.628'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.70%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '83.06'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +5.02%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.01%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.894'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.35%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.06'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.45%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0503'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.17%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.89'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.68%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.53%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.941.62'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '104.23'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +7.20%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.95'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.73%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '49.29'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.08%  '
$cell.Style = 'Normal'
